$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 251218
$ws.Range("B6").Value = "BIMEC 5"
$ws.Range("C6").Value = 21
$ws.Range("D6").Value = 96.90140845070422
$ws.Range("E6").Value = "2025-06-04 13:07:16"
$ws.Range("F6").Value = "2025-06-04 13:28:16"
$ws.Range("G6").Value = "2025-06-04 13:28:16"
$ws.Range("H6").Value = "2025-06-05 07:05:10"
$ws.Range("I6").Value = 6880
$ws.Range("J6").Value = "bobina"
$ws.Range("K6").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 76
$ws.Range("N6").Value = 39885
$ws.Range("O6").Value = "X"
$ws.Range("P6").Value = 39885
$ws.Range("Q6").Value = "2025-05-09 00:00:00"
$ws.Range("R6").Value = -0.2952660406828704
$ws.Range("S6").Value = 1

# Row 7
$ws.Range("A7").Value = 251895
$ws.Range("B7").Value = "BIMEC 5"
$ws.Range("C7").Value = 38
$ws.Range("D7").Value = 249.2112676056338
$ws.Range("E7").Value = "2025-06-05 07:05:10"
$ws.Range("F7").Value = "2025-06-05 07:43:10"
$ws.Range("G7").Value = "2025-06-05 07:43:10"
$ws.Range("H7").Value = "2025-06-05 11:52:23"
$ws.Range("I7").Value = 17694
$ws.Range("J7").Value = "bobina"
$ws.Range("K7").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$ws.Range("L7").Value = 10
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = "39891 (esterno)"
$ws.Range("O7").Value = "X"
$ws.Range("P7").Value = 39891
$ws.Range("Q7").Value = "2025-05-26 00:00:00"
$ws.Range("R7").Value = -10.49471830986111
$ws.Range("S7").Value = 4

# Row 8
$ws.Range("A8").Value = 251565
$ws.Range("B8").Value = "BIMEC 5"
$ws.Range("C8").Value = 31
$ws.Range("D8").Value = 176.7464788732394
$ws.Range("E8").Value = "2025-06-05 11:52:23"
$ws.Range("F8").Value = "2025-06-05 12:23:23"
$ws.Range("G8").Value = "2025-06-05 12:23:23"
$ws.Range("H8").Value = "2025-06-06 07:20:08"
$ws.Range("I8").Value = 12549
$ws.Range("J8").Value = "bobina"
$ws.Range("K8").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 39885
$ws.Range("O8").Value = "X"
$ws.Range("P8").Value = 39885
$ws.Range("Q8").Value = "2025-06-10 00:00:00"
$ws.Range("R8").Value = -1.30565336462963
$ws.Range("S8").Value = 1

# Row 9
$ws.Range("A9").Value = 251780
$ws.Range("B9").Value = "BIMEC 5"
$ws.Range("C9").Value = 17
$ws.Range("D9").Value = 342.2394366197183
$ws.Range("E9").Value = "2025-06-06 07:20:08"
$ws.Range("F9").Value = "2025-06-06 07:37:08"
$ws.Range("G9").Value = "2025-06-06 07:37:08"
$ws.Range("H9").Value = "2025-06-06 13:19:22"
$ws.Range("I9").Value = 24299
$ws.Range("J9").Value = "bobina"
$ws.Range("K9").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 70
$ws.Range("N9").Value = "39887 (esterno)"
$ws.Range("O9").Value = "X"
$ws.Range("P9").Value = 39887
$ws.Range("Q9").Value = "2025-06-18 00:00:00"
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 7

# Row 13
$ws.Range("A13").Value = 252084
$ws.Range("B13").Value = "R10"
$ws.Range("C13").Value = 20
$ws.Range("D13").Value = 641
$ws.Range("E13").Value = "2025-06-04 07:00:00"
$ws.Range("F13").Value = "2025-06-04 07:20:00"
$ws.Range("G13").Value = "2025-06-04 07:20:00"
$ws.Range("H13").Value = "2025-06-05 10:01:00"
$ws.Range("I13").Value = 39101
$ws.Range("J13").Value = "bobina"
$ws.Range("K13").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 70
$ws.Range("N13").Value = 39885
$ws.Range("O13").Value = "X"
$ws.Range("P13").Value = 39885
$ws.Range("Q13").Value = "2025-06-30 00:00:00"
$ws.Range("R13").Value = -0.4173611111111111
$ws.Range("S13").Value = 7

# Row 14
$ws.Range("A14").Value = 251984
$ws.Range("B14").Value = "R10"
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 338.327868852459
$ws.Range("E14").Value = "2025-06-05 10:01:00"
$ws.Range("F14").Value = "2025-06-05 10:26:00"
$ws.Range("G14").Value = "2025-06-05 10:26:00"
$ws.Range("H14").Value = "2025-06-06 08:04:19"
$ws.Range("I14").Value = 20638
$ws.Range("J14").Value = "bobina"
$ws.Range("K14").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L14").Value = 3
$ws.Range("M14").Value = 70
$ws.Range("N14").Value = 39874
$ws.Range("O14").Value = "X"
$ws.Range("P14").Value = 39874
$ws.Range("Q14").Value = "2025-06-10 00:00:00"
$ws.Range("R14").Value = -1.3363387978125
$ws.Range("S14").Value = 1

# Row 15
$ws.Range("A15").Value = 252282
$ws.Range("B15").Value = "R10"
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 44.88524590163934
$ws.Range("E15").Value = "2025-06-06 08:04:19"
$ws.Range("F15").Value = "2025-06-06 08:34:19"
$ws.Range("G15").Value = "2025-06-06 08:34:19"
$ws.Range("H15").Value = "2025-06-06 09:19:12"
$ws.Range("I15").Value = 2738
$ws.Range("J15").Value = "bobina"
$ws.Range("K15").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L15").Value = 5
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 39885
$ws.Range("O15").Value = "X"
$ws.Range("P15").Value = 39885
$ws.Range("Q15").Value = "2025-06-09 00:00:00"
$ws.Range("R15").Value = -1.388342440798611
$ws.Range("S15").Value = 1

# Row 16
$ws.Range("A16").Value = 251180
$ws.Range("B16").Value = "R10"
$ws.Range("C16").Value = 25
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = "2025-06-06 09:19:12"
$ws.Range("F16").Value = "2025-06-06 09:44:12"
$ws.Range("G16").Value = "2025-06-06 09:44:12"
$ws.Range("H16").Value = "2025-06-06 09:44:12"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "bobina"
$ws.Range("K16").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L16").Value = 4
$ws.Range("M16").Value = 70
$ws.Range("N16").Value = "39887 (esterno)"
$ws.Range("O16").Value = "X"
$ws.Range("P16").Value = 39887
$ws.Range("Q16").Value = "2025-05-20 00:00:00"
$ws.Range("R16").Value = -17.40570355190972
$ws.Range("S16").Value = 7

